$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("neg_reaction11")

$values = @(0,1,3,4,5,6,7,8,11,12,13,14,15,16,18,19,20,21,25,27,28,30,32,33,36,37,38,39,40,42,43,44,45,47,48,49,50,51,52,53,54,55,56,57,58,59,60,62,63,64,65,66,67,68,70,71,72,73,74,75)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $values[$i]
}
